# Update the PSSM score matrix (B2:K21) with the recomputed values from the
# supplemental-figures rerun. Baseline/"unseen" cells now carry
# -17.61049862271234 (was -19.7552934769669); residue-specific scores are
# replaced with their newly computed figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-17.61049862271234, 1.963061662625442, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234),
    @(-17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, 2.838195440519315, -17.61049862271234, -17.61049862271234),
    @(-17.61049862271234, 2.249741190590017, 2.1203247282511, -17.61049862271234, 3.452343992110698, -17.61049862271234, 1.56951005723913, -17.61049862271234, 4.321921246144655, -17.61049862271234),
    @(-17.61049862271234, 1.695137891213891, -17.61049862271234, -17.61049862271234, -17.61049862271234, 2.912750320876413, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234),
    @(-17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234),
    @(2.666808077040652, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234),
    @(-17.61049862271234, -17.61049862271234, -17.61049862271234, 1.767815286322248, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234),
    @(3.770795587628651, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234),
    @(-17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, 1.352355030844695, -17.61049862271234, 1.894640309495433),
    @(-17.61049862271234, -17.61049862271234, -17.61049862271234, 2.854714023689093, -17.61049862271234, 2.434931773262531, -17.61049862271234, -17.61049862271234, -17.61049862271234, 1.6483457103858),
    @(-17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234),
    @(-17.61049862271234, -17.61049862271234, -17.61049862271234, 2.353059320525733, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, 2.240678518889485),
    @(-17.61049862271234, -17.61049862271234, 1.285848564772182, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, 2.16544113643373),
    @(-17.61049862271234, -17.61049862271234, 1.279507311425621, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234),
    @(-17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234),
    @(-17.61049862271234, 2.108370339141272, 2.310065040697576, -17.61049862271234, -17.61049862271234, -17.61049862271234, 1.343348265381411, 1.276475008161301, -17.61049862271234, -17.61049862271234),
    @(-17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, -17.61049862271234, 0.8020422812548459, 0.8679905753614854, -17.61049862271234, -17.61049862271234),
    @(-17.61049862271234, -17.61049862271234, 1.502645789121098, -17.61049862271234, -17.61049862271234, -17.61049862271234, 1.631924226567218, 1.697273632264424, -17.61049862271234, -17.61049862271234),
    @(-17.61049862271234, 0.8467419625812694, 1.58184476611879, -17.61049862271234, 3.178526302340178, -17.61049862271234, 1.969040140413209, 1.487967885777805, -17.61049862271234, 1.97624039644663),
    @(-17.61049862271234, 0.9974854085989445, -17.61049862271234, 2.088099229564603, -17.61049862271234, 2.820077456061056, 2.520193845327887, -17.61049862271234, -17.61049862271234, -17.61049862271234)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}